$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.500856280326843
$ws.Range("B1").Value = 3.779188632965088
$ws.Range("C1").Value = 3.32874870300293
$ws.Range("D1").Value = 3.593692779541016
$ws.Range("E1").Value = 1.370967149734497
